$wb = $excel.ActiveWorkbook

# New trade rows to append (Trade #, Date, Time, Strategy, Side, Entry Price,
# Exit Price, Status, P&L %, P&L $, Capital After, Entry Slippage (bps),
# Exit Slippage (bps), Confidence, Entry Reason, Exit Reason, Duration (min))
$newTrades = @(
    @(53, "2026-02-17", "20:32:58", "MarketMaking", "DOWN", 0.3,  "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(54, "2026-02-17", "20:33:05", "MarketMaking", "UP",   0.79, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(55, "2026-02-17", "20:33:18", "MarketMaking", "DOWN", 0.2,  "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(56, "2026-02-17", "20:33:25", "MarketMaking", "DOWN", 0.22, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(57, "2026-02-17", "20:33:32", "MarketMaking", "UP",   0.67, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0),
    @(58, "2026-02-17", "20:33:39", "MarketMaking", "DOWN", 0.54, "OPEN", 0, 0, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", 0)
)

function Write-TradeRows($ws, $startRow) {
    $r = $startRow
    foreach ($trade in $newTrades) {
        $ws.Cells.Item($r, 1).Value = $trade[0]

        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $trade[1]

        $ws.Cells.Item($r, 3).NumberFormat = "@"
        $ws.Cells.Item($r, 3).Value = $trade[2]

        $ws.Cells.Item($r, 4).Value = $trade[3]
        $ws.Cells.Item($r, 5).Value = $trade[4]
        $ws.Cells.Item($r, 6).Value = $trade[5]
        $ws.Cells.Item($r, 7).Value = ""
        $ws.Cells.Item($r, 8).Value = $trade[6]
        $ws.Cells.Item($r, 9).Value = $trade[7]
        $ws.Cells.Item($r, 10).Value = $trade[8]
        $ws.Cells.Item($r, 11).Value = $trade[9]
        $ws.Cells.Item($r, 12).Value = $trade[10]
        $ws.Cells.Item($r, 13).Value = $trade[11]
        $ws.Cells.Item($r, 14).Value = $trade[12]
        $ws.Cells.Item($r, 15).Value = $trade[13]
        $ws.Cells.Item($r, 16).Value = ""
        $ws.Cells.Item($r, 17).Value = $trade[14]

        $r = $r + 1
    }
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Write-TradeRows $wsAllTrades 54

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Write-TradeRows $wsMarketMaking 21

Write-Output "Trading update applied: 2026-02-17 20:33:46"
